# Updated cryptos list on Fri Jul 21 17:28:28 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns with newly scraped
# coinranking.com figures, and reflects two rank swaps where the source
# rankings changed order (Solana/Dogecoin at rows 9-10, Quant/
# SynthetixNetwork at rows 48-49) so Coin/Link/Price/Volume all move
# together for those rows.
#
# Many Price values look like plain numbers (e.g. "0.9996"), but the
# workbook stores the Price column as text (so values like "29.840.09"
# round-trip exactly). To keep Excel from silently re-interpreting those
# cells as numbers, each numeric-looking target cell is temporarily
# switched to a text number format before its value is written, then the
# temporary format is cleared again so the cell's style is left exactly
# as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "29.849.07"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.891.39"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "0.7812"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").Value = "243.99"
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.3143"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "25.39"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.07313"
$ws.Range("E10").Value = "  +4.23%  "
$ws.Range("D11").Value = "0.08125"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "0.7666"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "5.467"
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("D14").Value = "1.876.82"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "93.16"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "6.205"
$ws.Range("D17").Value = "29.837.02"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "13.93"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "245.81"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("D21").Value = "0.9992"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "8.151"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("D23").Value = "2.128.07"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "0.1587"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").Value = "9.463"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").Value = "161.77"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").Value = "18.78"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "2.035"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").Value = "4.476"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D33").Value = "0.05602"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "0.7554"
$ws.Range("E36").Value = "  +2.62%  "
$ws.Range("D37").Value = "0.9978"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "2.639"
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("D39").Value = "0.01934"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").Value = "2.782"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "1.141.74"
$ws.Range("E41").Value = "  +10.95%  "
$ws.Range("D42").Value = "0.4451"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").Value = "73.89"
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("D44").Value = "5.973"
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("D45").Value = "0.8551"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").Value = "0.9996"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "1.902"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "101.88"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").Value = "3.106"
$ws.Range("E49").Value = "  +6.33%  "
$ws.Range("D50").Value = "9.811"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E51").Value = "  +0.97%  "

$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
